$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The sheet currently ends at row 105, which holds only a footnote
# ("*4/8..." shared string) in column B, with column A left blank.
# A new data row (date 2020-05-09 / serial 43960) needs to be inserted
# before that footnote, pushing the footnote row down to row 106.

# 1) Move the existing footnote row (row 105) down to row 106, copying
#    both its values and its formatting. Only columns A:B are populated
#    on that row, so only copy those two cells.
$ws.Range("A105:B105").Copy($ws.Range("A106:B106"))

# 2) Re-use the formatting pattern of the previous data row (104) for the
#    new data row 105.
$ws.Range("A104:E104").Copy($ws.Range("A105:E105"))

# 3) Fill in the new data values for row 105.
$ws.Range("A105").Value = 43960
$ws.Range("B105").Value = 378
$ws.Range("C105").Value = 35385
$ws.Range("D105").Value = 98
$ws.Range("E105").Value = 7232

# 4) Update the active selection to reflect the new last cell.
[void]$ws.Range("E106").Select()

# 5) Update the workbook's print area defined name to cover the new row.
for ($i = 1; $i -le $wb.Names.Count; $i++) {
    $n = $wb.Names.Item($i)
    if ($n.Name -like "*Print_Area*") {
        $n.RefersTo = "=相談件数!`$A`$1:`$E`$107"
    }
}
